$d = $word.ActiveDocument

# Find the paragraph that ends the "Break the problem apart:" answer for
# Problem 3 ("There are no constraints ...") so the three new paragraphs
# can be inserted right after it, ahead of the existing trailing blank
# paragraphs.
$search = $d.Content
$found = $search.Find.Execute(
    "There are no constraints that need to be considered.  The sub-goals are to have and effective method that will work with any given number.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$anchorIndex = $search.Paragraphs(1).Index

# 1) blank paragraph right after the anchor
$d.Paragraphs($anchorIndex).Range.InsertParagraphAfter()

# 2) "Identify potential solutions:" paragraph
$d.Paragraphs($anchorIndex + 1).Range.InsertParagraphAfter()
$d.Paragraphs($anchorIndex + 2).Range.Text = "Identify potential solutions:"

# 3) paragraph describing the two candidate solutions
$d.Paragraphs($anchorIndex + 2).Range.InsertParagraphAfter()
$d.Paragraphs($anchorIndex + 3).Range.Text = "First solution: create a massive chart the would display a column of each possible finger and then a row showing what numbers are counted for that finger.  Second solution:  create a mathematical formula that will help figure out which finger the count will stop on."
